$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1072470739083369, 0.583924513487991, 0.5944036395536997, 0.7709757710549013, 0.7806389925716273, 23)
    3  = @(0.573629067650584, 0.8880024064378264, 3.979434659499331, 1.994852039500507, 1.953538050233822, 23)
    4  = @(0.228737977167174, 1.392066633737142, 7.304711183963507, 2.702722920308981, 2.753551201305114, 23)
    5  = @(0.1376631994370348, 1.217887797378426, 7.521227355969471, 2.742485616365101, 2.80058721592678, 23)
    6  = @(0.256355043509169, 1.359960095653506, 7.005699073826143, 2.646828115655821, 2.693591441706439, 23)
    7  = @(0.09029628155329977, 1.369716599985631, 7.335564264182292, 2.708424683128976, 2.767756307023949, 23)
    8  = @(0.1622545047491004, 1.487177404830128, 8.202263652189526, 2.863959436198342, 2.92362286105626, 23)
    9  = @(0.08448834823307437, 1.390260216460968, 7.363914233767023, 2.713653300214864, 2.773296749308377, 23)
    10 = @(0.1037379453787874, 1.518653016633174, 7.970283297256951, 2.823169016771216, 2.88466950940459, 23)
    11 = @(-0.05287204651734589, 1.303960028208192, 7.512551892089196, 2.740903480987464, 2.801983048518609, 23)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
